# Update "想去人数" (F column) figures across sheets to match the refreshed
# data pull (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1714
$ws1.Range("F4").Value  = 123
$ws1.Range("F5").Value  = 513
$ws1.Range("F7").Value  = 1407
$ws1.Range("F8").Value  = 215
$ws1.Range("F9").Value  = 75
$ws1.Range("F10").Value = 113
$ws1.Range("F11").Value = 6027
$ws1.Range("F13").Value = 395
$ws1.Range("F15").Value = 4811
$ws1.Range("F16").Value = 22
$ws1.Range("F19").Value = 46
$ws1.Range("F20").Value = 349
$ws1.Range("F21").Value = 47
$ws1.Range("F25").Value = 3271
$ws1.Range("F26").Value = 135

# 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 58

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1714
$ws4.Range("F4").Value  = 123
$ws4.Range("F5").Value  = 58
$ws4.Range("F6").Value  = 513
$ws4.Range("F8").Value  = 1407
$ws4.Range("F9").Value  = 215
$ws4.Range("F10").Value = 75
$ws4.Range("F11").Value = 113
$ws4.Range("F12").Value = 6027
$ws4.Range("F14").Value = 395
$ws4.Range("F16").Value = 4811
$ws4.Range("F17").Value = 22
$ws4.Range("F20").Value = 46
$ws4.Range("F21").Value = 349
$ws4.Range("F22").Value = 47
$ws4.Range("F26").Value = 3272
$ws4.Range("F28").Value = 135
